# Update countries & provincias Spain
# - Re-rank three country rows (Barein/Afganistan, Bulgaria/Consejo Danes,
#   Togo/Jamaica each swap places with the row above them).
# - Refresh the daily case/death counters for the affected rows.
# - Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (column A) ---------------------------------------
$ws.Range("A49").Value = "Barein"
$ws.Range("A50").Value = "Afganistan"

$ws.Range("A84").Value = "Bulgaria"
$ws.Range("A85").Value = "Consejo Danes para los Refugiados"

$ws.Range("A150").Value = "Togo"
$ws.Range("A151").Value = "Jamaica"

# --- Refreshed statistics (columns B..H) ----------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3756573
$ws.Range("C4").Value = 61548
$ws.Range("D4").Value = 1708926
$ws.Range("E4").Value = 1905773
$ws.Range("G4").Value = 756
$ws.Range("H4").Value = 141874

# Row 5 - Brasil
$ws.Range("B5").Value = 2046328
$ws.Range("C5").Value = 31590
$ws.Range("E5").Value = 601702
$ws.Range("G5").Value = 1029
$ws.Range("H5").Value = 77851

# Row 10 - Chile
$ws.Range("E10").Value = 21378
$ws.Range("G10").Value = 98
$ws.Range("H10").Value = 8347

# Row 30 - Suecia
$ws.Range("B30").Value = 77281
$ws.Range("C30").Value = 152
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 5619

# Row 49 - now Barein
$ws.Range("B49").Value = 35473
$ws.Range("C49").Value = 389
$ws.Range("D49").Value = 31188
$ws.Range("E49").Value = 4161
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 124

# Row 50 - now Afganistan
$ws.Range("B50").Value = 35229
$ws.Range("C50").Value = 159
$ws.Range("D50").Value = 23151
$ws.Range("E50").Value = 10931
$ws.Range("G50").Value = 34
$ws.Range("H50").Value = 1147

# Row 84 - now Bulgaria
$ws.Range("B84").Value = 8442
$ws.Range("C84").Value = 298
$ws.Range("D84").Value = 4033
$ws.Range("E84").Value = 4112
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 297

# Row 85 - now Consejo Danes para los Refugiados
$ws.Range("B85").Value = 8249
$ws.Range("C85").Value = 50
$ws.Range("D85").Value = 4248
$ws.Range("E85").Value = 3808
$ws.Range("H85").Value = 193

# Row 87 - Estado de Palestina
$ws.Range("E87").Value = 6219
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 53

# Row 139 - Niger
$ws.Range("D139").Value = 1013
$ws.Range("E139").Value = 20

# Row 150 - now Togo
$ws.Range("B150").Value = 766
$ws.Range("C150").Value = 17
$ws.Range("D150").Value = 546
$ws.Range("E150").Value = 205
$ws.Range("H150").Value = 15

# Row 151 - now Jamaica
$ws.Range("B151").Value = 765
$ws.Range("C151").Value = 2
$ws.Range("D151").Value = 647
$ws.Range("E151").Value = 108
$ws.Range("H151").Value = 10

# Row 152 - Santo Tome y Principe
$ws.Range("B152").Value = 741
$ws.Range("C152").Value = 1
$ws.Range("E152").Value = 402

# --- Timestamp footer ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 23:54"
